# Anonymize "fedcore" -> "approach" and add header-row borders
# (underline style under the merged header cells), matching the
# target commit "update of results and scripts. Anonimyzed fedcore".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- Build the two new border styles once, on sheet1, then copy them
# --- to the equivalent cells on sheet2 so that the underlying style
# --- (cellXfs/borders) table stays minimal and gets reused instead of
# --- duplicated.

# C1 (sheet1): top + bottom thin border
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1Borders = $c1.Borders
$c1Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# D1 (sheet1): top + right + bottom thin border
$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1Borders = $d1.Borders
$d1Borders.Item(8).LineStyle = 1    # xlEdgeTop
$d1Borders.Item(10).LineStyle = 1   # xlEdgeRight
$d1Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# Reuse the freshly built styles for the matching cells on sheet2
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Remove the stray empty inline-string cell G5 on sheet2
$ws2.Range("G5").ClearContents()
